$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 13: update "Azami" (maximum) amounts
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Rows 24-25: clear the "600 TL" values in column E
$ws.Range("E24").Value = ""
$ws.Range("E25").Value = ""
